# Append a new row (row 26) of parsed packet data to each of the four
# "DE_xxx" log sheets, mirroring the row already present at row 25
# (same date-time style, same inline-string / numeric layout).

$wb = $excel.ActiveWorkbook

$newTimestamp = [double]"45812.43565972222"

$rowsData = @{
    1 = @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x78"
        E = "0x14"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 376
        I = 14
    }
    2 = @{
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x78"
        E = "0xe"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 376
        I = 14
    }
    3 = @{
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x82"
        E = "0x7"
        F = 130
        G = [double]"5.68631262647114e+23"
        H = 130
        I = 7
    }
    4 = @{
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x82"
        E = "0x3"
        F = 130
        G = [double]"9.85046333984776e+23"
        H = 130
        I = 3
    }
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $data = $rowsData[$i]

    $newRow = 26

    $ws.Cells.Item($newRow, 1).Value = $newTimestamp
    $ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($newRow - 1, 1).NumberFormat

    $ws.Cells.Item($newRow, 2).Value = $data.B
    $ws.Cells.Item($newRow, 3).Value = $data.C
    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 5).Value = $data.E
    $ws.Cells.Item($newRow, 6).Value = $data.F
    $ws.Cells.Item($newRow, 7).Value = $data.G
    $ws.Cells.Item($newRow, 8).Value = $data.H
    $ws.Cells.Item($newRow, 9).Value = $data.I
}
